$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6120544276223256
$ws.Range("C2").Value = 0.1582178864502168
$ws.Range("D2").Value = 0.01523005430277635
$ws.Range("F2").Value = 0.4835143216245967
$ws.Range("G2").Value = 0.002383005729340188
$ws.Range("M2").Value = 0.8931652751493857
$ws.Range("O2").Value = 1.544213758709077
$ws.Range("B3").Value = 0.5354954224460187
$ws.Range("C3").Value = 0.1438735245382929
$ws.Range("D3").Value = 0.01364120990714213
$ws.Range("F3").Value = 0.4778559183888476
$ws.Range("G3").Value = 0.002385830164463341
$ws.Range("M3").Value = 0.7947831649405543
$ws.Range("O3").Value = 1.539211341963835
$ws.Range("B4").Value = 0.4883364759589881
$ws.Range("C4").Value = 0.1350140981047332
$ws.Range("D4").Value = 0.01266017368185146
$ws.Range("F4").Value = 0.4748258946642494
$ws.Range("G4").Value = 0.00238765575187945
$ws.Range("M4").Value = 0.7349081208685533
$ws.Range("O4").Value = 1.537604654448927
$ws.Range("B5").Value = 0.4690823438641303
$ws.Range("C5").Value = 0.1313911699948562
$ws.Range("D5").Value = 0.01225905561926055
$ws.Range("F5").Value = 0.4737024137879118
$ws.Range("G5").Value = 0.002388422744542693
$ws.Range("M5").Value = 0.7106360203322595
$ws.Range("O5").Value = 1.537316707355927
$ws.Range("B6").Value = 0.4658830538158156
$ws.Range("C6").Value = 0.1307888344805122
$ws.Range("D6").Value = 0.01219237064919554
$ws.Range("F6").Value = 0.473522568877442
$ws.Range("G6").Value = 0.002388551497554459
$ws.Range("M6").Value = 0.7066131732506733
$ws.Range("O6").Value = 1.53729100189156
$ws.Range("B7").Value = 0.4880769535794229
$ws.Range("C7").Value = 0.1349652886500508
$ws.Range("D7").Value = 0.01265476942693056
$ws.Range("F7").Value = 0.4748102930174483
$ws.Range("G7").Value = 0.002387666002497192
$ws.Range("M7").Value = 0.7345802709710085
$ws.Range("O7").Value = 1.537599288008721
$ws.Range("B8").Value = 0.5856892369392881
$ws.Range("C8").Value = 0.1532829743556761
$ws.Range("D8").Value = 0.01468338142039727
$ws.Range("F8").Value = 0.4814708339673217
$ws.Range("G8").Value = 0.002383960675176058
$ws.Range("M8").Value = 0.8591288566878603
$ws.Range("O8").Value = 1.54218396431412
$ws.Range("B9").Value = 0.7758446251762621
$ws.Range("C9").Value = 0.1887751886500553
$ws.Range("D9").Value = 0.01861644025704834
$ws.Range("F9").Value = 0.4980781488172781
$ws.Range("G9").Value = 0.002377416132534079
$ws.Range("M9").Value = 1.107889796336323
$ws.Range("O9").Value = 1.562870050734006
$ws.Range("B10").Value = 0.9147147020558464
$ws.Range("C10").Value = 0.2145700613854444
$ws.Range("D10").Value = 0.02147677101582701
$ws.Range("F10").Value = 0.5124724733267669
$ws.Range("G10").Value = 0.002373042968827365
$ws.Range("M10").Value = 1.293851259506113
$ws.Range("O10").Value = 1.58530332332279
$ws.Range("B11").Value = 0.9776951630773851
$ws.Range("C11").Value = 0.2262397142641532
$ws.Range("D11").Value = 0.02277128312764631
$ws.Range("F11").Value = 0.5195038344906067
$ws.Range("G11").Value = 0.002371146959983345
$ws.Range("M11").Value = 1.379243064672607
$ws.Range("O11").Value = 1.597102804373236
$ws.Range("B12").Value = 1.001515178961256
$ws.Range("C12").Value = 0.2306490380301227
$ws.Range("D12").Value = 0.0232604862612078
$ws.Range("F12").Value = 0.5222364374263719
$ws.Range("G12").Value = 0.002370442338494412
$ws.Range("M12").Value = 1.411701292022158
$ws.Range("O12").Value = 1.601802019915738
$ws.Range("B13").Value = 0.996386442808614
$ws.Range("C13").Value = 0.229699850499884
$ws.Range("D13").Value = 0.02315517277432377
$ws.Range("F13").Value = 0.5216448025116449
$ws.Range("G13").Value = 0.00237059349832902
$ws.Range("M13").Value = 1.404705273834793
$ws.Range("O13").Value = 1.600779658113822
$ws.Range("B14").Value = 0.9796554472461025
$ws.Range("C14").Value = 0.2266026692857395
$ws.Range("D14").Value = 0.02281155045417904
$ws.Range("F14").Value = 0.5197272421455921
$ws.Range("G14").Value = 0.002371088723079779
$ws.Range("M14").Value = 1.381910929178645
$ws.Range("O14").Value = 1.597484772025126
$ws.Range("B15").Value = 0.9694033618080766
$ws.Range("C15").Value = 0.2247042775674686
$ws.Range("D15").Value = 0.02260094032845927
$ws.Range("F15").Value = 0.5185618086332937
$ws.Range("G15").Value = 0.002371393799736665
$ws.Range("M15").Value = 1.367964878409282
$ws.Range("O15").Value = 1.59549669448586
$ws.Range("B16").Value = 0.9105947691350025
$ws.Range("C16").Value = 0.2138060855659774
$ws.Range("D16").Value = 0.02139203387061883
$ws.Range("F16").Value = 0.5120227237712243
$ws.Range("G16").Value = 0.002373168750356581
$ws.Range("M16").Value = 1.2882873866522
$ws.Range("O16").Value = 1.584564431732758
$ws.Range("B17").Value = 0.8744671712182139
$ws.Range("C17").Value = 0.2071035499609764
$ws.Range("D17").Value = 0.02064867154293637
$ws.Range("F17").Value = 0.5081353355735416
$ws.Range("G17").Value = 0.002374281488101579
$ws.Range("M17").Value = 1.239617298629568
$ws.Range("O17").Value = 1.578267451672502
$ws.Range("B18").Value = 0.8536695378502941
$ws.Range("C18").Value = 0.2032423775570464
$ws.Range("D18").Value = 0.02022048417346411
$ws.Range("F18").Value = 0.5059448755206688
$ws.Range("G18").Value = 0.002374930297764657
$ws.Range("M18").Value = 1.211698102470493
$ws.Range("O18").Value = 1.574795571378644
$ws.Range("B19").Value = 0.846624777189561
$ws.Range("C19").Value = 0.2019340249685797
$ws.Range("D19").Value = 0.02007540122051665
$ws.Range("F19").Value = 0.5052110170570501
$ws.Range("G19").Value = 0.002375151485796819
$ws.Range("M19").Value = 1.202257718487033
$ws.Range("O19").Value = 1.573645758380025
$ws.Range("B20").Value = 0.878314890231195
$ws.Range("C20").Value = 0.207817675764403
$ws.Range("D20").Value = 0.02072786872330568
$ws.Range("F20").Value = 0.5085444459828921
$ws.Range("G20").Value = 0.002374162126027035
$ws.Range("M20").Value = 1.244790542441166
$ws.Range("O20").Value = 1.578922241513823
$ws.Range("B21").Value = 0.984570555175253
$ws.Range("C21").Value = 0.22751265375922
$ws.Range("D21").Value = 0.02291250814445078
$ws.Range("F21").Value = 0.520288573071312
$ws.Range("G21").Value = 0.002370942901725299
$ws.Range("M21").Value = 1.388602802617868
$ws.Range("O21").Value = 1.598446276679965
$ws.Range("B22").Value = 1.053843118737063
$ws.Range("C22").Value = 0.2403276014283904
$ws.Range("D22").Value = 0.02433444535879659
$ws.Range("F22").Value = 0.5283721535040939
$ws.Range("G22").Value = 0.002368916769226141
$ws.Range("M22").Value = 1.483308399076364
$ws.Range("O22").Value = 1.612553657595384
$ws.Range("B23").Value = 1.01688730476792
$ws.Range("C23").Value = 0.2334933686619536
$ws.Range("D23").Value = 0.02357607981613086
$ws.Range("F23").Value = 0.5240202933962763
$ws.Range("G23").Value = 0.002369991057158988
$ws.Range("M23").Value = 1.432694165157642
$ws.Range("O23").Value = 1.604900434438974
$ws.Range("B24").Value = 0.8765754216706227
$ws.Range("C24").Value = 0.2074948437367539
$ws.Range("D24").Value = 0.02069206622371667
$ws.Range("F24").Value = 0.5083593487474616
$ws.Range("G24").Value = 0.002374216061407409
$ws.Range("M24").Value = 1.242451526725588
$ws.Range("O24").Value = 1.578625749118828
$ws.Range("B25").Value = 0.724544712778993
$ws.Range("C25").Value = 0.1792215712837049
$ws.Range("D25").Value = 0.01755745605343151
$ws.Range("F25").Value = 0.4932025813634056
$ws.Range("G25").Value = 0.002379109849696222
$ws.Range("M25").Value = 1.040064438344004
$ws.Range("O25").Value = 1.5560109503337
